# Auto update on 2026-01-08 11:56:14
# Applies updated numeric values to the kp_data export on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 94
$ws.Range("O3").Value = 94
$ws.Range("R3").Value = 0.41

# Row 4
$ws.Range("F4").Value = 84
$ws.Range("N4").Value = 338
$ws.Range("P4").Value = 114
$ws.Range("Q4").Value = 2.96

# Row 6
$ws.Range("J6").Value = 438
$ws.Range("N6").Value = 438
$ws.Range("P6").Value = 86
$ws.Range("Q6").Value = 5.09

# Row 7
$ws.Range("D7").Value = 162
$ws.Range("J7").Value = 230
$ws.Range("N7").Value = 416
$ws.Range("P7").Value = 98
$ws.Range("Q7").Value = 4.24
$ws.Range("R7").Value = 0.1

# Row 9
$ws.Range("F9").Value = 193
$ws.Range("N9").Value = 193
$ws.Range("Q9").Value = 2.44

# Row 10
$ws.Range("F10").Value = 247
$ws.Range("G10").Value = 88
$ws.Range("N10").Value = 247
$ws.Range("O10").Value = 88
$ws.Range("P10").Value = 78
$ws.Range("Q10").Value = 3.17

# Row 11
$ws.Range("F11").Value = 49
$ws.Range("N11").Value = 49
$ws.Range("P11").Value = 24
$ws.Range("Q11").Value = 2.04

# Row 12
$ws.Range("F12").Value = 103
$ws.Range("N12").Value = 340
$ws.Range("P12").Value = 107

# Row 13
$ws.Range("F13").Value = 136
$ws.Range("N13").Value = 136
$ws.Range("P13").Value = 51
$ws.Range("Q13").Value = 2.67
$ws.Range("R13").Value = 0.1

# Row 14
$ws.Range("F14").Value = 331
$ws.Range("G14").Value = 55
$ws.Range("N14").Value = 331
$ws.Range("O14").Value = 55
$ws.Range("P14").Value = 100
$ws.Range("Q14").Value = 3.31
$ws.Range("R14").Value = 0.17

# Row 15
$ws.Range("F15").Value = 52
$ws.Range("G15").Value = 9
$ws.Range("N15").Value = 95
$ws.Range("O15").Value = 9
$ws.Range("P15").Value = 32
$ws.Range("Q15").Value = 2.97

# Row 16
$ws.Range("F16").Value = 72
$ws.Range("G16").Value = 8
$ws.Range("N16").Value = 75
$ws.Range("O16").Value = 8
$ws.Range("P16").Value = 45
$ws.Range("Q16").Value = 1.67
$ws.Range("R16").Value = 0.11

# Row 17
$ws.Range("F17").Value = 29
$ws.Range("N17").Value = 29
$ws.Range("Q17").Value = 3.62
